# Refresh the crypto price/volume snapshot (GitHub Actions symbol-list update).
# Price (col D) and Volume/1h (col E) are stored as *text* in the sheet, not
# numbers, so numeric-looking values are written with a leading apostrophe to
# force Excel to keep them as text (preserving exact formatting such as
# trailing/leading zeros, e.g. "31.10" or "0.00006015").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'295.54"
$ws.Range("E2").Value = "'-2.44%"
$ws.Range("D3").Value = "'31.10"
$ws.Range("E3").Value = "'-2.47%"
$ws.Range("D4").Value = "'5.114"
$ws.Range("E4").Value = "'-2.21%"
$ws.Range("D5").Value = "'0.07362"
$ws.Range("E5").Value = "'1.49%"
$ws.Range("E6").Value = "'-1.21%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.744"
$ws.Range("E7").Value = "'-0.16%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.649"
$ws.Range("E8").Value = "'12.93%"
$ws.Range("D9").Value = "'0.9172"
$ws.Range("E9").Value = "'0.75%"
$ws.Range("D10").Value = "'0.1673"
$ws.Range("E10").Value = "'-0.39%"
$ws.Range("D11").Value = "'0.07230"
$ws.Range("E11").Value = "'-2.42%"
$ws.Range("D12").Value = "'0.07913"
$ws.Range("E12").Value = "'-0.82%"
$ws.Range("D13").Value = "'0.02982"
$ws.Range("E13").Value = "'0.00%"
$ws.Range("D14").Value = "'0.09914"
$ws.Range("E14").Value = "'-0.16%"
$ws.Range("D15").Value = "'0.001492"
$ws.Range("E15").Value = "'-1.10%"
$ws.Range("D16").Value = "'0.006157"
$ws.Range("E16").Value = "'-6.11%"
$ws.Range("D17").Value = "'3.447"
$ws.Range("E17").Value = "'-1.06%"
$ws.Range("E18").Value = "'0.12%"
$ws.Range("E19").Value = "'-1.73%"
$ws.Range("D20").Value = "'0.1347"
$ws.Range("E20").Value = "'1.85%"
$ws.Range("D21").Value = "'4.555"
$ws.Range("E21").Value = "'6.11%"
$ws.Range("D22").Value = "'0.04618"
$ws.Range("E22").Value = "'1.21%"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'-1.15%"
$ws.Range("D25").Value = "'0.004422"
$ws.Range("E25").Value = "'0.20%"
$ws.Range("E26").Value = "'-1.32%"
$ws.Range("D27").Value = "'0.0001873"
$ws.Range("E27").Value = "'6.96%"
$ws.Range("D39").Value = "'0.01679"
$ws.Range("E39").Value = "'1.04%"
$ws.Range("D40").Value = "'0.04421"
$ws.Range("E40").Value = "'-1.35%"
$ws.Range("D41").Value = "'0.007085"
$ws.Range("E41").Value = "'1.62%"
$ws.Range("D42").Value = "'0.1328"
$ws.Range("E42").Value = "'-1.33%"
$ws.Range("D43").Value = "'0.002104"
$ws.Range("E43").Value = "'-12.90%"
$ws.Range("D44").Value = "'0.01101"
$ws.Range("E44").Value = "'-14.17%"
$ws.Range("D45").Value = "'0.00006015"
$ws.Range("E45").Value = "'-1.31%"
$ws.Range("D47").Value = "'0.01022"
$ws.Range("E47").Value = "'-21.79%"
